$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 56, shifting existing rows 56-79 down to 57-80
$ws.Rows.Item(56).Insert()

# Populate the newly inserted row 56 with the new weekly record
$ws.Range("A56").Value = 10
$ws.Range("B56").Value = "Vega Modelo de Temuco"
$ws.Range("C56").Value = "La Araucanía"
$ws.Range("D56").Value = 44466
$ws.Range("E56").Value = 9
$ws.Range("F56").Value = 100112012
$ws.Range("G56").Value = "Espinaca"
$ws.Range("H56").Value = "Sin especificar"
$ws.Range("I56").Value = "Primera"
$ws.Range("J56").Value = 40
$ws.Range("K56").Value = 10000
$ws.Range("L56").Value = 10000
$ws.Range("M56").Value = 10000
$ws.Range("N56").Value = "$/docena de atados"
$ws.Range("O56").Value = "Región de La Araucanía"
$ws.Range("P56").Value = 3333
$ws.Range("Q56").Value = 3
$ws.Range("R56").Value = "Hortaliza"
